$wb = $excel.ActiveWorkbook

# --- Insert a new sheet "2022-Q4" right after "总计" (before current "2022-Q3") ---
$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ3Old = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($sheetQ3Old)
$newSheet.Name = "2022-Q4"

# --- Populate the new "2022-Q4" sheet with fund holdings data ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'015032"
$newSheet.Range("C2").Value = "中融医药消费混合A"
$newSheet.Range("D2").Value = "'0.52"
$newSheet.Range("E2").Value = "'92.98"
$newSheet.Range("F2").Value = "'4.41"
$newSheet.Range("G2").Value = "'0.0229"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'015033"
$newSheet.Range("C3").Value = "中融医药消费混合C"
$newSheet.Range("D3").Value = "'0.04"
$newSheet.Range("E3").Value = "'92.98"
$newSheet.Range("F3").Value = "'4.41"
$newSheet.Range("G3").Value = "'0.0018"
$newSheet.Range("H3").Value = 4

# --- Match the bold/bordered "header row" + "index column" formatting used
# throughout the workbook by copying it over from the sibling "2022-Q3" sheet
# (re-fetched by name now that the sheet collection has shifted) ---
$sheetQ3Ref = $wb.Worksheets.Item("2022-Q3")
$sheetQ3Ref.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$sheetQ3Ref.Range("A2:A3").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- Update the "总计" (summary) sheet: shift existing quarter rows down and
# insert the new 2022-Q4 row at the top ---
$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.02

$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2022-Q3"
$sheetTotal.Range("C3").Value = 7
$sheetTotal.Range("D3").Value = 0.27

$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("B4").Value = "2022-Q2"
$sheetTotal.Range("C4").Value = 6
$sheetTotal.Range("D4").Value = 0.33

$sheetTotal.Range("A5").Value = 3
$sheetTotal.Range("B5").Value = "2022-Q1"
$sheetTotal.Range("C5").Value = 8
$sheetTotal.Range("D5").Value = 0.73

# Match formatting for the newly added index cell in column A (row 5)
$sheetTotal.Range("A4").Copy()
$sheetTotal.Range("A5").PasteSpecial(-4122)
$sheetTotal.Range("A5").Value = 3
